$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'67.623.98"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -1.22%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.789.97"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.95%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.03%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'595.23"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -0.03%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'166.84"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.07%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'3.788.04"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.97%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.09%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.519"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.28%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.49%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'6.36"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -2.04%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.449"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.19%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.74%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'36.43"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.81%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'4.422.10"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.96%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'3.797.89"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +1.49%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  +3.21%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'67.575.29"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -1.22%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.27%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.21%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'10.16"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -5.21%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'456.86"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -2.36%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.697"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.20%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'0.0000156"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +7.79%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  -1.03%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'11.94"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.91%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  -2.47%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'10.08"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -0.36%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.00%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  +0.51%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'7.29"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +0.07%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = "'EthereumClassic"
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = "'29.86"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.00%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = "'ImmutableX"
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = "'2.20"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.68%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'9.24"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -0.13%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'0.999"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.38%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'3.741.25"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.88%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.100"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -1.03%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.90%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.74%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.998"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.46%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'5.77"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -0.66%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.999"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.00%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'Arweave"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'45.76"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +5.84%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'USDe"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'1.00"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -0.01%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  -1.70%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'47.17"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +3.14%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'148.85"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +1.59%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'8.33"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -3.03%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  -4.65%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'389.22"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.55%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  +1.57%  "
$ws.Range('E51').Style = 'Normal'
